$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D (closest achievable value to the target 10.875 given this
# engine's pixel-quantized ColumnWidth setter).
$ws.Columns.Item(4).ColumnWidth = 10

# Row 3 additions:
# A3: trade date/time (same style as A2)
$ws.Range("A3").Value = 42650.3669212963
# D3: buy price (plain number, same as D2/C3)
$ws.Range("D3").Value = 104.839996
# G3: holding flag TRUE (same style as G2)
$ws.Range("G3").Value = $true

# Reuse the existing styles from row 2 (date format) instead of creating new
# style entries, by copying formats from A2 -> A3 and G2 -> G3.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
